$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @{ B = 73177.4875;          C = 5444.512500000001;  D = 65000; E = 2732.975000000006 }
    3  = @{ B = 69811.6565;          C = 5387.3435;          D = 45000; E = 19424.31299999999 }
    4  = @{ B = 64462.7775;          C = 5357.2225;          D = 22500; E = 36605.55499999999 }
    5  = @{ B = 64419.8455;          C = 5289.154500000001;  D = 22500; E = 36630.69100000001 }
    6  = @{ B = 63786.1715;          C = 5361.8285;          D = 22500; E = 35924.34299999999 }
    7  = @{ B = 65583.3985;          C = 5442.6015;          D = 22500; E = 37640.797 }
    8  = @{ B = 63552.469;           C = 5816.530999999999;  D = 22500; E = 35235.93799999999 }
    9  = @{ B = 76403.3885;          C = 7034.6115;          D = 22500; E = 46868.777 }
    10 = @{ B = 59351.3105;          C = 8596.6895;          D = 65000; E = -14245.379 }
    11 = @{ B = 99192.12700000001;   C = 13942.873;          D = 65000; E = 20249.25400000002 }
    12 = @{ B = 99767.6905;          C = 15769.3095;         D = 65000; E = 18998.38099999999 }
    13 = @{ B = 100144.771;          C = 15209.229;          D = 65000; E = 19935.54199999999 }
    14 = @{ B = 101117.205;          C = 15159.795;          D = 65000; E = 20957.41 }
    15 = @{ B = 108724.5555;         C = 15449.4445;         D = 65000; E = 28275.111 }
    16 = @{ B = 103763.567;          C = 15588.433;          D = 65000; E = 23175.13399999999 }
    17 = @{ B = 101211.3415;         C = 16020.6585;         D = 65000; E = 20190.68299999999 }
    18 = @{ B = 93839.61749999999;   C = 16262.3825;         D = 65000; E = 12577.23499999999 }
    19 = @{ B = 88583.391;           C = 15862.609;          D = 65000; E = 7720.782000000007 }
    20 = @{ B = 89024.58199999999;   C = 14887.418;          D = 65000; E = 9137.16399999999 }
    21 = @{ B = 87063.96249999999;   C = 13341.0375;         D = 65000; E = 8722.924999999988 }
    22 = @{ B = 85910.32799999999;   C = 11837.672;          D = 65000; E = 9072.655999999988 }
    23 = @{ B = 84684.008;           C = 9490.991999999998;  D = 65000; E = 10193.016 }
    24 = @{ B = 81619.20600000001;   C = 6891.794;           D = 65000; E = 9727.412000000011 }
    25 = @{ B = 52348.0285;          C = 5770.9715;          D = 65000; E = -18422.943 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("B$row").Value = $vals.B
    $ws.Range("C$row").Value = $vals.C
    $ws.Range("D$row").Value = $vals.D
    $ws.Range("E$row").Value = $vals.E
}
